$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "a"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Insira o Nome do Sistema"
